$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F7").Value = 2819
$sheet1.Range("F8").Value = 1678
$sheet1.Range("F9").Value = 1818
$sheet1.Range("F12").Value = 744
$sheet1.Range("F13").Value = 895
$sheet1.Range("F14").Value = 167
$sheet1.Range("F15").Value = 376
$sheet1.Range("F16").Value = 1114
$sheet1.Range("F18").Value = 50
$sheet1.Range("F20").Value = 6605
$sheet1.Range("F21").Value = 253
$sheet1.Range("F22").Value = 1522
$sheet1.Range("F23").Value = 159
$sheet1.Range("F24").Value = 178
$sheet1.Range("F25").Value = 158
$sheet1.Range("F26").Value = 308
$sheet1.Range("F27").Value = 267
$sheet1.Range("F28").Value = 67
$sheet1.Range("F29").Value = 1100
$sheet1.Range("F30").Value = 906
$sheet1.Range("F32").Value = 90
$sheet1.Range("F34").Value = 478
$sheet1.Range("F35").Value = 1465
$sheet1.Range("F36").Value = 159
$sheet1.Range("F37").Value = 141
$sheet1.Range("F38").Value = 218
$sheet1.Range("F39").Value = 15
$sheet1.Range("F40").Value = 143
$sheet1.Range("F41").Value = 201
$sheet1.Range("F42").Value = 158

$sheet2 = $wb.Worksheets.Item("演出")
$sheet2.Range("F8").Value = 17

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F10").Value = 2819
$sheet4.Range("F11").Value = 1678
$sheet4.Range("F12").Value = 1818
$sheet4.Range("F15").Value = 744
$sheet4.Range("F17").Value = 895
$sheet4.Range("F18").Value = 167
$sheet4.Range("F19").Value = 376
$sheet4.Range("F20").Value = 1114
$sheet4.Range("F21").Value = 50
$sheet4.Range("F23").Value = 6605
$sheet4.Range("F24").Value = 253
$sheet4.Range("F25").Value = 1522
$sheet4.Range("F27").Value = 159
$sheet4.Range("F28").Value = 178
$sheet4.Range("F29").Value = 158
$sheet4.Range("F30").Value = 308
$sheet4.Range("F31").Value = 267
$sheet4.Range("F32").Value = 67
$sheet4.Range("F33").Value = 1100
$sheet4.Range("F34").Value = 906
$sheet4.Range("F36").Value = 90
$sheet4.Range("F38").Value = 478
$sheet4.Range("F39").Value = 1465
$sheet4.Range("F40").Value = 159
$sheet4.Range("F41").Value = 141
$sheet4.Range("F42").Value = 218
$sheet4.Range("F43").Value = 15
$sheet4.Range("F44").Value = 143
$sheet4.Range("F45").Value = 201
$sheet4.Range("F47").Value = 17
$sheet4.Range("F49").Value = 158
